$d = $word.ActiveDocument

# Locate the target paragraph: "If Nu/s<1 then there are no single mutants ... increasing the mutation rate ..."
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*then there are no single mutants*") {
        $targetPara = $cand
        break
    }
}

$pStart = $targetPara.Range.Start
$pEnd = $targetPara.Range.End
$pText = $d.Range($pStart, $pEnd).Text

# --- Part 1: relocate the "_GoBack" bookmark away from before
#     "then there are no single mutants..." (it moves further along the
#     paragraph, right before the replacement text inserted in Part 2). The
#     visible text around the old bookmark position is unchanged.
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# --- Part 2: replace the highlighted sentence's ending ---
# "increasing the mutation rate of individuals with fitness below 1 will not affect the process."
# becomes two new, unhighlighted runs:
#   "increasing the mutation rate of individuals with fitness below 1 will "
#   "have a much smaller effect than if single mutants were abundant"
# followed by a highlighted "." (new run, reusing the yellow highlight).
# The relocated "_GoBack" bookmark is placed right before this new text.
$oldSentence = "increasing the mutation rate of individuals with fitness below 1 will not affect the process."
$sIdx = $pText.IndexOf($oldSentence)
$sStart = $pStart + $sIdx
$sEnd = $sStart + $oldSentence.Length

$oldRange = $d.Range($sStart, $sEnd)
$oldRange.Delete()

# Relocate the bookmark to right before the new replacement text.
$d.Bookmarks.Add("_GoBack", $d.Range($sStart, $sStart))

# Insert the new plain (non-highlighted) text.
$newPlain = "increasing the mutation rate of individuals with fitness below 1 will have a much smaller effect than if single mutants were abundant"
$insPlain = $d.Range($sStart, $sStart)
$insPlain.InsertAfter($newPlain)
$plainEnd = $sStart + $newPlain.Length

# Insert the highlighted final period right after the new plain text.
$insPeriod = $d.Range($plainEnd, $plainEnd)
$insPeriod.InsertAfter(".")
$periodRange = $d.Range($plainEnd, $plainEnd + 1)
$periodRange.HighlightColorIndex = "yellow"

$finalCheck = $d.Range($targetPara.Range.Start, $targetPara.Range.End)
Write-Output $finalCheck.Text
